$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item(1)
$ws.Range("F2").Value = 158
$ws.Range("F3").Value = 931
$ws.Range("F4").Value = 1092
$ws.Range("F5").Value = 1553
$ws.Range("F8").Value = 12577
$ws.Range("F9").Value = 2215
$ws.Range("F11").Value = 277
$ws.Range("F12").Value = 21676
$ws.Range("F14").Value = 1253
$ws.Range("F15").Value = 243
$ws.Range("F17").Value = 809
$ws.Range("F18").Value = 684
$ws.Range("F19").Value = 324
$ws.Range("F20").Value = 2941
$ws.Range("F21").Value = 790
$ws.Range("F22").Value = 4508
$ws.Range("F23").Value = 1164
$ws.Range("F24").Value = 889
$ws.Range("F29").Value = 1107
$ws.Range("F30").Value = 61
$ws.Range("F31").Value = 124
$ws.Range("F32").Value = 286
$ws.Range("F35").Value = 41
$ws.Range("F36").Value = 26
$ws.Range("F37").Value = 4512
$ws.Range("F39").Value = 4622
$ws.Range("F40").Value = 5587
$ws.Range("F42").Value = 133
$ws.Range("F43").Value = 96
$ws.Range("F44").Value = 182
$ws.Range("F45").Value = 373
$ws.Range("F47").Value = 52
$ws.Range("F48").Value = 4123
$ws.Range("F49").Value = 152

$ws = $wb.Worksheets.Item(2)
$ws.Range("F2").Value = 4172
$ws.Range("F4").Value = 71
$ws.Range("F5").Value = 105
$ws.Range("F12").Value = 1059

$ws = $wb.Worksheets.Item(3)
$ws.Range("F3").Value = 491
$ws.Range("F4").Value = 109
$ws.Range("F5").Value = 20

$ws = $wb.Worksheets.Item(4)
$ws.Range("F2").Value = 491
$ws.Range("F3").Value = 158
$ws.Range("F4").Value = 931
$ws.Range("F5").Value = 1092
$ws.Range("F6").Value = 1553
$ws.Range("F9").Value = 12577
$ws.Range("F10").Value = 2215
$ws.Range("F12").Value = 277
$ws.Range("F13").Value = 1253
$ws.Range("F14").Value = 243
$ws.Range("F16").Value = 809
$ws.Range("F17").Value = 684
$ws.Range("F18").Value = 324
$ws.Range("F19").Value = 2941
$ws.Range("F20").Value = 790
$ws.Range("F21").Value = 4508
$ws.Range("F22").Value = 4508
$ws.Range("F23").Value = 1164
$ws.Range("F24").Value = 20
$ws.Range("F25").Value = 105
$ws.Range("F31").Value = 1107
$ws.Range("F33").Value = 124
$ws.Range("F35").Value = 286
$ws.Range("F38").Value = 41
$ws.Range("F39").Value = 4622
$ws.Range("F40").Value = 133
$ws.Range("F41").Value = 182
$ws.Range("F46").Value = 4123
